# Update column F ("dSF") values on Sheet1 to reflect the repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    5  = 1
    6  = -1
    7  = 1
    8  = 4
    9  = 3
    10 = -2
    11 = -2
    12 = 2
    13 = -1
    14 = 1
    15 = -3
    16 = -3
    17 = -1
    18 = -2
    19 = 5
    21 = 3
    22 = 3
    23 = 7
    24 = -1
    25 = -5
    26 = -3
    27 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
